# Adding the required columns required for the field report impact taxonomy
#
# The taxonomy sheet already has an "orgtypengo / IFRC / GO-Appeal" block
# (rows 31-33). This adds a parallel "orgtypengo / IFRC / GO-Field Reports"
# block (new rows 34-65), one row per field-report variable name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new rows with the same look & feel (fonts/fills) as the existing
# GO-Appeal rows directly above them, columns A:D (the only columns used by
# this block).
$ws.Range("A33:D33").Copy()
$ws.Range("A34:D65").PasteSpecial(-4122)

$fields = @(
  "num_injured","num_dead","num_missing","num_affected","num_displaced",
  "num_assisted","num_localstaff","num_volunteers","num_expats_delegates",
  "num_potentially_affected","num_highest_risk",
  "gov_num_injured","gov_num_dead","gov_num_missing","gov_num_affected",
  "gov_num_displaced","gov_num_assisted","health_min_num_assisted",
  "gov_num_potentially_affected","gov_num_highest_risk",
  "other_num_injured","other_num_dead","other_num_missing","other_num_affected",
  "other_num_displaced","other_num_assisted","other_num_potentially_affected",
  "other_num_highest_risk",
  "dref_amount","appeal_amount","imminent_dref_amount","forecast_based_action_amount"
)

$row = 34
foreach ($f in $fields) {
  $ws.Range("A$row").Value = "orgtypengo"
  $ws.Range("B$row").Value = "IFRC"
  $ws.Range("C$row").Value = "GO-Field Reports"
  $ws.Range("D$row").Value = $f
  $row = $row + 1
}

# Reset the view: scroll back so the top-left visible cell is A1, and select
# the first empty row beneath the newly-added data (A66:B66) - matching the
# cursor position left behind after typing in the last new row.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A66:B66").Select()
